$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 new columns before column D; this shifts existing D:K data to F:M
$ws.Range("D:E").Insert()

# New columns inherit format from the left (col C); copy the number format from
# column F (the original column D, now shifted) so D:E match the data columns,
# restricted to the rows that actually hold table data
$ws.Range("F7:F35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("F7:F35").Copy()
$ws.Range("E7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns (D, E) with the new quarter figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 238300
$ws.Range("E8").Value = 207700
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 1300
$ws.Range("E15").Value = 1100
$ws.Range("D17").Value = 175100
$ws.Range("E17").Value = 154800
$ws.Range("D18").Value = 63200
$ws.Range("E18").Value = 52900
$ws.Range("D20").Value = 5100
$ws.Range("E20").Value = 1600
$ws.Range("D21").Value = 69600
$ws.Range("E21").Value = 55700
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 68200
$ws.Range("E23").Value = 54500
$ws.Range("D24").Value = 12200
$ws.Range("E24").Value = 9600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 56000
$ws.Range("E26").Value = 44900
$ws.Range("D27").Value = 39200
$ws.Range("E27").Value = 32400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -5100
$ws.Range("E32").Value = -1600
$ws.Range("D33").Value = 39200
$ws.Range("E33").Value = 32400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 39200
$ws.Range("E35").Value = 32400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 261100
$ws.Range("E41").Value = 86900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 68600
$ws.Range("E43").Value = 81500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 15700
$ws.Range("E45").Value = 14600
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 143900
$ws.Range("E47").Value = 204700
$ws.Range("D48").Value = 12700
$ws.Range("E48").Value = 12000
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 412300
$ws.Range("E52").Value = 423400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 914400
$ws.Range("E54").Value = 823000
$ws.Range("D57").Value = 19800
$ws.Range("E57").Value = 16900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 197700
$ws.Range("E59").Value = 142300
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 510000
$ws.Range("E66").Value = 452500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -237800
$ws.Range("E72").Value = -250600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 404400
$ws.Range("E76").Value = 370500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 39200
$ws.Range("E81").Value = 32400
$ws.Range("D83").Value = 1300
$ws.Range("E83").Value = 1100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 158300
$ws.Range("E89").Value = 162600
$ws.Range("D91").Value = -2000
$ws.Range("E91").Value = -1800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 63000
$ws.Range("E94").Value = -87600
$ws.Range("D96").Value = -27700
$ws.Range("E96").Value = -119000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -47400
$ws.Range("E100").Value = -119600
$ws.Range("D101").Value = 200
$ws.Range("E101").Value = -1200
$ws.Range("D102").Value = 174100
$ws.Range("E102").Value = -45800

# Row 91 (Capital Expenditures) received fully revised figures, not a simple shift
$ws.Range("D91").Value = -2000
$ws.Range("E91").Value = -1800
$ws.Range("F91").Value = -2100
$ws.Range("G91").Value = -1100
$ws.Range("H91").Value = -1100
$ws.Range("I91").Value = -1500
$ws.Range("J91").Value = -1200
$ws.Range("K91").Value = -1800
$ws.Range("L91").Value = -100
$ws.Range("M91").Value = 1800

$ws.Range("D7:E102").Columns.AutoFit()
